$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.849.72'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '1.709.48'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.91'
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3967'
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4116'
$ws.Range("E8").Value = '  +2.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.528'
$ws.Range("E9").Value = '  +0.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.003'
$ws.Range("E10").Value = '  +0.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.58'
$ws.Range("E11").Value = '  +2.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08997'
$ws.Range("E12").Value = '  +2.71%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.669'
$ws.Range("E13").Value = '  +6.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.37'
$ws.Range("E14").Value = '  +4.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.200'
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001365'
$ws.Range("E16").Value = '  +3.98%  '
$ws.Range("D17").Value = '1.710.33'
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '100.32'
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07149'
$ws.Range("E19").Value = '  +1.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.19'
$ws.Range("E20").Value = '  +2.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.493'
$ws.Range("E21").Value = '  +6.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.008'
$ws.Range("E22").Value = '  +0.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.57'
$ws.Range("E23").Value = '  +2.41%  '
$ws.Range("D24").Value = '24.876.01'
$ws.Range("E24").Value = '  +0.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.093'
$ws.Range("E25").Value = '  -1.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.352'
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.08'
$ws.Range("E27").Value = '  +1.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.319'
$ws.Range("E28").Value = '  +24.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.75'
$ws.Range("E29").Value = '  +2.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '139.52'
$ws.Range("E30").Value = '  +1.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.241'
$ws.Range("E31").Value = '  +0.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.884'
$ws.Range("E32").Value = '  +9.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09089'
$ws.Range("E33").Value = '  +5.79%  '
$ws.Range("E34").Value = '  +0.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03033'
$ws.Range("E35").Value = '  +10.85%  '
$ws.Range("E36").Value = '  +3.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.18'
$ws.Range("E37").Value = '  -3.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.964'
$ws.Range("E38").Value = '  +2.34%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8331'
$ws.Range("E39").Value = '  +8.84%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '14.66'
$ws.Range("E40").Value = '  +1.63%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09340'
$ws.Range("E41").Value = '  +2.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.492'
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.66'
$ws.Range("E43").Value = '  +6.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7530'
$ws.Range("E44").Value = '  +5.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.661'
$ws.Range("E45").Value = '  +1.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.280'
$ws.Range("E46").Value = '  +1.43%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.003'
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("B48").Value = 'Flow'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.352'
$ws.Range("E48").Value = '  +1.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.16'
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '94.41'
$ws.Range("E50").Value = '  +5.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08087'
$ws.Range("E51").Value = '  +1.07%  '
